# fix(excel): table columns init
#
# The two tabs had their header rows filled with the wrong (duplicated)
# label and ended up associated with the wrong sheet name. Rename the
# tabs to match their real content and rewrite the header rows with the
# correct, distinct column names.

$wb = $excel.ActiveWorkbook

# --- First tab: was "CONTACTS" (FIRSTNAME x6) -> becomes "USERS" ---------
$wsUsers = $wb.Worksheets.Item(1)
$wsContactsTmp = $wb.Worksheets.Item(2)

# Swap names via a temporary name so we never collide with the existing
# sheet name while the rename is in flight.
$wsContactsTmp.Name = "__TMP_SWAP__"
$wsUsers.Name = "USERS"

# Clear the old, wider header row (A1:F1) before writing the new one.
$wsUsers.Range("A1:F1").ClearContents()

$wsUsers.Range("A1").Value = "USERNAME"
$wsUsers.Range("B1").Value = "PASSWORD"
$wsUsers.Range("C1").Value = "CREATE_AT"

# --- Second tab: was "USERS" (USERNAME x3) -> becomes "CONTACTS" ---------
$wsContacts = $wsContactsTmp
$wsContacts.Name = "CONTACTS"

$wsContacts.Range("A1").Value = "FIRSTNAME"
$wsContacts.Range("B1").Value = "MIDDLENAME"
$wsContacts.Range("C1").Value = "LASTNAME"
$wsContacts.Range("D1").Value = "PHONE"
$wsContacts.Range("E1").Value = "EMAIL"
$wsContacts.Range("F1").Value = "COMMENT"
$wsContacts.Range("G1").Value = "CREATED_AT"
